# Daily attendance processing - 2025-12-09 10:31:48
#
# The "Recorded By" column (G) stores a comma-separated list of the
# accounts that recorded/updated each attendance session. This pass
# rotates each multi-value list left by one position (the first
# recorder moves to the end of the list), leaving single-value cells
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $v = $cell.Value2

    if ($v -ne $null -and $v -is [string] -and $v.Contains(",")) {
        $parts = $v -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value2 = $rotated -join ", "
        }
    }
}
